$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet holds one weekly snapshot per row for "Haba" at Mapocho Venta
# Directa de Santiago. The weekly refresh shifts each existing week's data
# down into a new row ordering and appends a brand-new week's record.

# Row 2
$ws.Range("D2").Value = 44467
$ws.Range("J2").Value = 35
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 480

# Row 3
$ws.Range("D3").Value = 44425
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("P3").Value = 560

# Row 4
$ws.Range("D4").Value = 44340
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 15000
$ws.Range("L4").Value = 15000
$ws.Range("M4").Value = 15000
$ws.Range("P4").Value = 600

# Row 5
$ws.Range("D5").Value = 44449
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 16000
$ws.Range("M5").Value = 16000
$ws.Range("P5").Value = 640

# Row 6
$ws.Range("D6").Value = 44453
$ws.Range("J6").Value = 55
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14455
$ws.Range("P6").Value = 578

# Row 7
$ws.Range("D7").Value = 44418
$ws.Range("J7").Value = 12
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("P7").Value = 600

# Row 8
$ws.Range("D8").Value = 44432
$ws.Range("O8").Value = "Provincia del Elquí"

# Row 9
$ws.Range("D9").Value = 44446
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 13000
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 520

# Row 10
$ws.Range("D10").Value = 45113
$ws.Range("J10").Value = 8
$ws.Range("K10").Value = 17000
$ws.Range("L10").Value = 17000
$ws.Range("M10").Value = 17000
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 680

# Row 11
$ws.Range("D11").Value = 44435
$ws.Range("J11").Value = 15
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("M11").Value = 14000
$ws.Range("P11").Value = 560

# Row 12
$ws.Range("D12").Value = 44435
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 14000
$ws.Range("O12").Value = "Provincia del Elquí"
$ws.Range("P12").Value = 560

# Row 13
$ws.Range("D13").Value = 44421
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 15000
$ws.Range("P13").Value = 600

# New row 14 - append a new weekly record
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C14").Value = "Metropolitana"
$ws.Range("D14").Value = 44376
$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112026
$ws.Range("G14").Value = "Haba"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 15
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 12000
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Provincia de Limarí"
$ws.Range("P14").Value = 480
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
